# Update gh-pages output data ("想去人数" / interest counters in column F)
# across the relevant worksheets, matching the commit's regenerated snapshot.

$wb = $excel.ActiveWorkbook

# -------- Sheet "展览" (Exhibitions) --------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 3318
$wsExpo.Range("F6").Value  = 4812
$wsExpo.Range("F7").Value  = 465
$wsExpo.Range("F8").Value  = 293
$wsExpo.Range("F11").Value = 282
$wsExpo.Range("F14").Value = 656
$wsExpo.Range("F20").Value = 4747
$wsExpo.Range("F24").Value = 5885
$wsExpo.Range("F29").Value = 4418
$wsExpo.Range("F31").Value = 93
$wsExpo.Range("F33").Value = 855
$wsExpo.Range("F34").Value = 73
$wsExpo.Range("F36").Value = 783
$wsExpo.Range("F37").Value = 826

# -------- Sheet "本地生活" (Local life) --------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 1091

# -------- Sheet "全部类型" (All types - aggregate of the other sheets) --------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1091
$wsAll.Range("F8").Value  = 3318
$wsAll.Range("F10").Value = 4812
$wsAll.Range("F11").Value = 465
$wsAll.Range("F12").Value = 293
$wsAll.Range("F15").Value = 282
$wsAll.Range("F18").Value = 656
$wsAll.Range("F25").Value = 4747
$wsAll.Range("F29").Value = 5885
$wsAll.Range("F34").Value = 4418
$wsAll.Range("F37").Value = 93
$wsAll.Range("F39").Value = 855
$wsAll.Range("F40").Value = 73
$wsAll.Range("F42").Value = 783
$wsAll.Range("F43").Value = 826
